# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" right before the "总计" (totals) sheet
#    and populate it with the per-fund holdings detail for 2022-Q1.
# 2) Insert a new leading data row into "总计" summarizing 2022-Q1
#    (holdings count = 21, market value = 16.98), pushing the existing
#    quarters down by one row and renumbering the index column.

$wb = $excel.ActiveWorkbook

# NOTE: worksheet handles returned by Worksheets.Item(...) are positional —
# inserting/adding a sheet shifts what later index a held variable resolves
# to. So re-fetch sheets by name immediately after any Worksheets.Add() call
# instead of reusing a handle obtained beforehand.

# --- 1. Create & position the new "2022-Q1" sheet -------------------------
$q1 = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$q1.Name = "2022-Q1"

$templateSheet = $wb.Worksheets.Item("2021-Q4")

# Match header / index-column formatting used by the other quarter sheets.
$templateSheet.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$templateSheet.Range("A2:A22").Copy()
$q1.Range("A2:A22").PasteSpecial(-4122)

$q1.Cells.Item(1,2).Value = "基金代码"
$q1.Cells.Item(1,3).Value = "基金名称"
$q1.Cells.Item(1,4).Value = "基金规模"
$q1.Cells.Item(1,5).Value = "股票总仓位"
$q1.Cells.Item(1,6).Value = "仓位占比"
$q1.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1.Cells.Item(1,8).Value = "仓位排名"

# Columns B (fund code, has leading zeros) and D:G (text-formatted numbers)
# must stay text so they match the other quarter sheets and don't lose
# leading zeros / turn into shared numeric values.
$q1.Range("B2:B22").NumberFormat = "@"
$q1.Range("D2:G22").NumberFormat = "@"

$rows = @(
    @("002943","广发多因子灵活配置混合","209.03","94.20","4.04","8.4448",4),
    @("004702","南方金融主题灵活配置混合","30.93","89.81","7.64","2.3631",2),
    @("011401","汇添富成长精选混合A","42.03","87.99","3.33","1.3996",9),
    @("550008","信诚优胜精选混合","27.24","79.82","3.16","0.8608",5),
    @("006392","中信保诚创新成长灵活配置混合","25.22","81.30","2.85","0.7188",5),
    @("161838","银华创业板两年定期开放混合","10.44","95.40","5.72","0.5972",8),
    @("550002","信诚精萃成长混合","22.09","80.84","2.60","0.5743",9),
    @("000480","东方红新动力灵活配置混合","15.38","72.90","3.03","0.4660",8),
    @("001564","东方红京东大数据灵活配置混合","11.40","69.58","3.11","0.3545",6),
    @("159851","华宝中证金融科技主题ETF","3.16","98.58","6.25","0.1975",3),
    @("000219","博时裕益灵活配置混合","3.90","91.22","4.91","0.1915",4),
    @("159966","华夏创业板低波蓝筹ETF","7.65","98.96","2.38","0.1821",10),
    @("001167","金鹰科技创新股票","4.03","94.55","4.47","0.1801",6),
    @("210009","金鹰核心资源混合","3.86","94.96","4.66","0.1799",7),
    @("010908","大成沪深300指数增强A","3.00","86.52","2.90","0.0870",10),
    @("011402","汇添富成长精选混合C","1.65","87.99","3.33","0.0549",9),
    @("516100","华夏中证金融科技主题交易型开放式指数证券投资基金","0.68","96.91","6.18","0.0420",3),
    @("010909","大成沪深300指数增强C","1.41","86.52","2.90","0.0409",10),
    @("168701","合煦智远国证香蜜湖金融科技指数(LOF)A","0.90","93.15","2.12","0.0191",9),
    @("002802","广发东财大数据精选灵活配置混合","0.41","55.13","4.01","0.0164",1),
    @("168702","合煦智远国证香蜜湖金融科技指数(LOF)C","0.22","93.15","2.12","0.0047",9)
)

$r = 2
foreach ($row in $rows) {
    $q1.Cells.Item($r,1).Value = ($r - 2)
    $q1.Cells.Item($r,2).Value = $row[0]
    $q1.Cells.Item($r,3).Value = $row[1]
    $q1.Cells.Item($r,4).Value = $row[2]
    $q1.Cells.Item($r,5).Value = $row[3]
    $q1.Cells.Item($r,6).Value = $row[4]
    $q1.Cells.Item($r,7).Value = $row[5]
    $q1.Cells.Item($r,8).Value = $row[6]
    $r = $r + 1
}

# --- 2. Insert the 2022-Q1 summary row into "总计" -------------------------
# Re-fetch "总计" by name (its position shifted after the sheet was added).
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q1"
$totalSheet.Cells.Item(2,3).Value = 21
$totalSheet.Cells.Item(2,4).Value = 16.98

# The inserted row copies formatting from the row above (the header); match
# the plain index-column style ("s=2") used by the rest of the data rows.
$templateSheet.Range("A2").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Cells.Item(2,1).Value = 0

# Renumber the (now shifted) index column for the remaining quarters.
for ($row = 3; $row -le 7; $row++) {
    $totalSheet.Cells.Item($row,1).Value = $row - 2
}
